# Boost MOSFET Dissipation Calcs.xlsx - working-tree fix / SubBoards + FlatPack BOM push
#
# Updates the Input Variables (row 4), the SubBoards weight table (B30:B34),
# and the FlatPack BOM parameters (row 38) - plus turns I38 into a live
# formula that just mirrors F38 instead of a stale hard-coded number.
# All the dependent Pcond / Psw / Fsw formulas recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Input Variables (row 4) ------------------------------------------------
$ws.Range("D4").Value2 = 0.04
$ws.Range("E4").Value2 = 5.5
$ws.Range("I4").Value2 = 54

# --- SubBoards weights (B30:B34) -------------------------------------------
$ws.Range("B30").Value2 = 1192
$ws.Range("B31").Value2 = 1216
$ws.Range("B32").Value2 = 1252
$ws.Range("B33").Value2 = 842
$ws.Range("B34").Value2 = 624

# --- FlatPack BOM parameters (row 38) ---------------------------------------
$ws.Range("C38").Value2 = 15
$ws.Range("F38").Value2 = 600000
# I38 used to be a stray hard-coded value; now just track F38 directly.
$ws.Range("I38").Formula = "=F38"

# --- Restore the view: scroll back to the top and move the selection -------
$excel.Goto($ws.Range("E5"), $true)
